# Lattice multiplication exercises update
# Replaces the contents of every cell in the single table with a new
# set of multiplication problems, per the target revision.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Each entry: row, col, top ("AA x BB"), digits row ("  D    D"),
# left column digit 1, left column digit 2.
$cellsData = @(
    @{ Row = 1; Col = 1; Top = "64 x 86"; Digits = "  8    6"; Left1 = "6"; Left2 = "4" },
    @{ Row = 1; Col = 2; Top = "17 x 91"; Digits = "  9    1"; Left1 = "1"; Left2 = "7" },
    @{ Row = 1; Col = 3; Top = "98 x 86"; Digits = "  8    6"; Left1 = "9"; Left2 = "8" },

    @{ Row = 2; Col = 1; Top = "80 x 40"; Digits = "  4    0"; Left1 = "8"; Left2 = "0" },
    @{ Row = 2; Col = 2; Top = "64 x 29"; Digits = "  2    9"; Left1 = "6"; Left2 = "4" },
    @{ Row = 2; Col = 3; Top = "64 x 82"; Digits = "  8    2"; Left1 = "6"; Left2 = "4" },

    @{ Row = 3; Col = 1; Top = "69 x 36"; Digits = "  3    6"; Left1 = "6"; Left2 = "9" },
    @{ Row = 3; Col = 2; Top = "82 x 17"; Digits = "  1    7"; Left1 = "8"; Left2 = "2" },
    @{ Row = 3; Col = 3; Top = "26 x 79"; Digits = "  7    9"; Left1 = "2"; Left2 = "6" },

    @{ Row = 4; Col = 1; Top = "96 x 70"; Digits = "  7    0"; Left1 = "9"; Left2 = "6" },
    @{ Row = 4; Col = 2; Top = "34 x 85"; Digits = "  8    5"; Left1 = "3"; Left2 = "4" },
    @{ Row = 4; Col = 3; Top = "49 x 45"; Digits = "  4    5"; Left1 = "4"; Left2 = "9" },

    @{ Row = 5; Col = 1; Top = "56 x 76"; Digits = "  7    6"; Left1 = "5"; Left2 = "6" },
    @{ Row = 5; Col = 2; Top = "98 x 70"; Digits = "  7    0"; Left1 = "9"; Left2 = "8" },
    @{ Row = 5; Col = 3; Top = "45 x 60"; Digits = "  6    0"; Left1 = "4"; Left2 = "5" }
)

foreach ($cellData in $cellsData) {
    $cell = $t.Cell($cellData.Row, $cellData.Col)

    $left1Text = "$($cellData.Left1)|    |"
    $left2Text = "$($cellData.Left2)|    |"

    $runXml = "<w:r>" +
        "<w:rPr><w:sz w:val=`"32`"/></w:rPr>" +
        "<w:t>$($cellData.Top)</w:t>" +
        "<w:br/>" +
        "<w:t xml:space=`"preserve`">$($cellData.Digits)</w:t>" +
        "<w:br/>" +
        "<w:t xml:space=`"preserve`">  ----</w:t>" +
        "<w:br/>" +
        "<w:t>$left1Text</w:t>" +
        "<w:br/>" +
        "<w:t>$left2Text</w:t>" +
        "</w:r>"

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        "<w:body><w:p>$runXml</w:p></w:body>" +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $cell.Range.InsertXML($xml)
}
